# Regenerate orders with updated distance/size codes.
# The workbook stores condition/filename strings built from distance codes
# (D64/D80/D51) and size codes (S25/S20/S30). This run renumbers them to
# D69/D86/D55 and S31 (S30 -> S31), leaving S25/S20 untouched.
#
# Applied as whole-workbook substring replacements (LookAt:=xlPart is the
# Excel default for Range.Replace), which correctly rewrites every
# occurrence across the Condition, Filename_Left, Filename_Right, Distance,
# and Size columns in one pass, matching how the values are composed
# (e.g. "Face10_D64_S25" -> "Face10_D69_S25", "Fixation_D64_l.png" ->
# "Fixation_D69_l.png", "D64" -> "D69", "S30" -> "S31").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$xlPart = 2
$xlByRows = 1

$r1 = $rng.Replace("D64", "D69", $xlPart, $xlByRows, $true, $false, $false)
$r2 = $rng.Replace("D80", "D86", $xlPart, $xlByRows, $true, $false, $false)
$r3 = $rng.Replace("D51", "D55", $xlPart, $xlByRows, $true, $false, $false)
$r4 = $rng.Replace("S30", "S31", $xlPart, $xlByRows, $true, $false, $false)
